$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.092.69'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.879.11'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.04%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5039'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3965'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08218'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.66%  '
$ws.Range('D13').Value = '1.883.15'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.304'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.207'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.82%  '
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001086'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06476'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').Value = '30.087.69'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.845'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.158'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.29%  '
$ws.Range('D26').Value = '2.091.97'
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.241'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.081'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  -2.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.942'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.693'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.294'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06371'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.87%  '
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.174'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.514'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.55%  '
$ws.Range('E41').Value = '  -3.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.218'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.99%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.16'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5913'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.12%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.100'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.635'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.56%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.42'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.209'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.28%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.71%  '
